$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 27.07
$ws.Range("F2").Value = 19.41
$ws.Range("K2").Value = 60.7
$ws.Range("N2").Value = 49.16024380385575

# Row 3
$ws.Range("D3").Value = 51.75
$ws.Range("E3").Value = 57.2
$ws.Range("F3").Value = 4.97
$ws.Range("K3").Value = 56.7
$ws.Range("N3").Value = 49.16024380385575

# Row 4
$ws.Range("D4").Value = 28.08
$ws.Range("F4").Value = 9.800000000000001
$ws.Range("K4").Value = 56.1
$ws.Range("N4").Value = 49.16024380385575

# Row 5
$ws.Range("D5").Value = 308.79
$ws.Range("E5").Value = 52.9
$ws.Range("F5").Value = 0.07000000000000001
$ws.Range("G5").Value = 50
$ws.Range("K5").Value = 56.1
$ws.Range("N5").Value = 49.16024380385575
